$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the whole "Meta description" paragraph (the document's
#    2nd paragraph), but keep its content so it can be reused below
#    (it supplies both the new bold heading run and the body copy
#    that replaces the old "Create a feature image" prompt text).
# ------------------------------------------------------------------
$metaPara = $d.Paragraphs(2)
$metaPara.Range.Cut()

# ------------------------------------------------------------------
# 2) Paste the cut paragraph back in, right before the very last
#    paragraph of the document (the former "Create a feature image"
#    prompt paragraph), creating a brand-new paragraph there.
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$insertionPoint = $lastPara.Range.Start
$insertRange = $d.Range($insertionPoint, $insertionPoint)
$insertRange.Paste()

# ------------------------------------------------------------------
# 3) The pasted paragraph currently reads:
#       "Meta description: Play African Elephant for free and
#        experience an exciting safari-themed slot game. Dive into
#        the vivid savannahs of Africa and win big."
#    with "Meta description" bold and the remainder plain text.
#    Turn it into the new bold heading paragraph:
#       "Play African Elephant for Free - Exciting Casino Game"
# ------------------------------------------------------------------
$n = $d.Paragraphs.Count
$newHeadingPara = $d.Paragraphs($n - 1)

$headingScope = $newHeadingPara.Range
$headingScope.Find.ClearFormatting()
$headingScope.Find.Execute("Meta description", $true, $false, $false, $false, $false, $true, 1, $false, "Play African Elephant for Free - Exciting Casino Game", 1) | Out-Null

$oldTail = ": Play African Elephant for free and experience an exciting safari-themed slot game. Dive into the vivid savannahs of Africa and win big."
$headingScope2 = $newHeadingPara.Range
$headingScope2.Find.ClearFormatting()
$headingScope2.Find.Execute($oldTail, $true, $false, $false, $false, $false, $true, 1, $false, "", 1) | Out-Null

# ------------------------------------------------------------------
# 4) Replace the italic "Create a feature image..." prompt text
#    (now the very last paragraph) with the meta-description copy.
# ------------------------------------------------------------------
$oldPrompt = "Create a feature image for African Elephant that showcases a happy Maya warrior with glasses on a safari tour, surrounded by the game's iconic animals such as a cheetah, buffalo, and of course, an African elephant. The image should be in a cartoon style that's colorful and vibrant, with elements of the savanna in the background, like trees and grass. Place the Maya warrior front and center, with a big smile on his face and a camera in hand, ready to capture the animals' beauty. The image should also include the game's logo and some of the game's symbols, like the elephant and the playing cards. Make sure the image exudes excitement and adventure to entice players to try out the game."
$newMeta = "Play African Elephant for free and experience an exciting safari-themed slot game. Dive into the vivid savannahs of Africa and win big."

$finalPara = $d.Paragraphs($d.Paragraphs.Count)
$finalScope = $finalPara.Range
$finalScope.Find.ClearFormatting()
$finalScope.Find.Execute($oldPrompt, $true, $false, $false, $false, $false, $true, 1, $false, $newMeta, 1) | Out-Null

Write-Output "edit complete"
